$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6682722
$ws.Range("F23").Value = 'Talaea El Geish'
$ws.Range("G23").Value = 'Ceramica Cleopatra'
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 'D'
$ws.Range("K23").Value = 2.5
$ws.Range("L23").Value = 2.7
$ws.Range("M23").Value = 2.9
$ws.Range("N23").Value = 3.1
$ws.Range("O23").Value = 2.6
$ws.Range("P23").Value = 2.5
$ws.Range("Q23").Value = 0.25
$ws.Range("R23").Value = 1.7
$ws.Range("S23").Value = 2.1
$ws.Range("T23").Value = 1.75
$ws.Range("U23").Value = 1.75
$ws.Range("V23").Value = 2.05
$ws.Range("X23").Value = 1.6
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = 0.35
$ws.Range("AA23").Value = -0.5
$ws.Range("AC23").Value = 1.05
$ws.Range("B24").Value = 6682721
$ws.Range("F24").Value = 'Al Moqawloon Al Arab'
$ws.Range("G24").Value = 'El Daklyeh'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 'A'
$ws.Range("K24").Value = 2.05
$ws.Range("L24").Value = 2.9
$ws.Range("M24").Value = 3.5
$ws.Range("N24").Value = 2.3
$ws.Range("O24").Value = 2.8
$ws.Range("P24").Value = 3.1
$ws.Range("Q24").Value = -0.25
$ws.Range("R24").Value = 1.975
$ws.Range("S24").Value = 1.825
$ws.Range("T24").Value = 2.25
$ws.Range("U24").Value = 2
$ws.Range("V24").Value = 1.8
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = 2.1
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0.825
$ws.Range("AC24").Value = 0.8
$ws.Range("B31").Value = 6853128
$ws.Range("F31").Value = 'Ghazl El Mahallah'
$ws.Range("G31").Value = 'National Bank'
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 'D'
$ws.Range("K31").Value = 2.875
$ws.Range("L31").Value = 2.75
$ws.Range("M31").Value = 2.45
$ws.Range("N31").Value = 2.55
$ws.Range("O31").Value = 2.75
$ws.Range("P31").Value = 2.75
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = 1.775
$ws.Range("S31").Value = 2.025
$ws.Range("T31").Value = 2.25
$ws.Range("U31").Value = 1.9
$ws.Range("V31").Value = 1.9
$ws.Range("W31").Value = -1
$ws.Range("X31").Value = 1.75
$ws.Range("Z31").Value = 0
$ws.Range("AA31").Value = -0
$ws.Range("AC31").Value = 0.8999999999999999
$ws.Range("B32").Value = 6853140
$ws.Range("F32").Value = 'Talaea El Geish'
$ws.Range("G32").Value = 'Aswan FC'
$ws.Range("H32").Value = 1
$ws.Range("J32").Value = 'H'
$ws.Range("K32").Value = 2.05
$ws.Range("L32").Value = 2.8
$ws.Range("M32").Value = 3.75
$ws.Range("N32").Value = 2
$ws.Range("O32").Value = 2.8
$ws.Range("P32").Value = 4
$ws.Range("Q32").Value = -0.5
$ws.Range("R32").Value = 2.05
$ws.Range("S32").Value = 1.75
$ws.Range("T32").Value = 2
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 1.8
$ws.Range("W32").Value = 1
$ws.Range("X32").Value = -1
$ws.Range("Z32").Value = 1.05
$ws.Range("AA32").Value = -1
$ws.Range("AC32").Value = 0.8
$ws.Range("B33").Value = 6853129
$ws.Range("F33").Value = 'Ceramica Cleopatra'
$ws.Range("G33").Value = 'Pyramids FC'
$ws.Range("H33").Value = 2
$ws.Range("J33").Value = 'H'
$ws.Range("K33").Value = 3.4
$ws.Range("L33").Value = 2.9
$ws.Range("M33").Value = 2.1
$ws.Range("N33").Value = 3
$ws.Range("O33").Value = 3.1
$ws.Range("P33").Value = 2.3
$ws.Range("Q33").Value = 0.25
$ws.Range("R33").Value = 1.775
$ws.Range("S33").Value = 2.025
$ws.Range("T33").Value = 2.25
$ws.Range("U33").Value = 2.025
$ws.Range("V33").Value = 1.775
$ws.Range("W33").Value = 2
$ws.Range("Y33").Value = -1
$ws.Range("Z33").Value = 0.7749999999999999
$ws.Range("AA33").Value = -1
$ws.Range("AB33").Value = 1.025
$ws.Range("AC33").Value = -1
$ws.Range("B34").Value = 6853139
$ws.Range("F34").Value = 'El Daklyeh'
$ws.Range("G34").Value = 'Ismaily SC'
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 'A'
$ws.Range("K34").Value = 2.7
$ws.Range("L34").Value = 2.7
$ws.Range("M34").Value = 2.7
$ws.Range("N34").Value = 2.8
$ws.Range("O34").Value = 2.75
$ws.Range("P34").Value = 2.55
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = 1.975
$ws.Range("S34").Value = 1.825
$ws.Range("T34").Value = 2
$ws.Range("U34").Value = 1.775
$ws.Range("V34").Value = 2.025
$ws.Range("W34").Value = -1
$ws.Range("Y34").Value = 1.55
$ws.Range("Z34").Value = -1
$ws.Range("AA34").Value = 0.825
$ws.Range("AB34").Value = -1
$ws.Range("AC34").Value = 1.025
$ws.Range("B38").Value = 6853133
$ws.Range("F38").Value = 'Pyramids FC'
$ws.Range("G38").Value = 'Talaea El Geish'
$ws.Range("H38").Value = 4
$ws.Range("I38").Value = 2
$ws.Range("K38").Value = 1.727
$ws.Range("M38").Value = 4.5
$ws.Range("N38").Value = 1.6
$ws.Range("O38").Value = 3.3
$ws.Range("P38").Value = 5.5
$ws.Range("Q38").Value = -1
$ws.Range("R38").Value = 1.975
$ws.Range("S38").Value = 1.825
$ws.Range("T38").Value = 2.5
$ws.Range("U38").Value = 1.9
$ws.Range("V38").Value = 1.9
$ws.Range("W38").Value = 0.6000000000000001
$ws.Range("Z38").Value = 0.9750000000000001
$ws.Range("AB38").Value = 0.8999999999999999
$ws.Range("B39").Value = 6853132
$ws.Range("F39").Value = 'Haras El Hedoud'
$ws.Range("G39").Value = 'Ghazl El Mahallah'
$ws.Range("H39").Value = 2
$ws.Range("I39").Value = 1
$ws.Range("K39").Value = 5.5
$ws.Range("M39").Value = 1.6
$ws.Range("N39").Value = 4.5
$ws.Range("O39").Value = 3.4
$ws.Range("P39").Value = 1.666
$ws.Range("Q39").Value = 0.75
$ws.Range("R39").Value = 1.825
$ws.Range("S39").Value = 1.975
$ws.Range("T39").Value = 2.25
$ws.Range("U39").Value = 1.8
$ws.Range("V39").Value = 2
$ws.Range("W39").Value = 3.5
$ws.Range("Z39").Value = 0.825
$ws.Range("AB39").Value = 0.8
$ws.Range("B40").Value = 6853143
$ws.Range("F40").Value = 'National Bank'
$ws.Range("G40").Value = 'Coca Cola FC'
$ws.Range("H40").Value = 1
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 'H'
$ws.Range("K40").Value = 2.75
$ws.Range("L40").Value = 2.9
$ws.Range("M40").Value = 2.5
$ws.Range("N40").Value = 2.5
$ws.Range("O40").Value = 2.8
$ws.Range("P40").Value = 2.9
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = 1.725
$ws.Range("S40").Value = 2.075
$ws.Range("W40").Value = 1.5
$ws.Range("Y40").Value = -1
$ws.Range("Z40").Value = 0.7250000000000001
$ws.Range("AA40").Value = -1
$ws.Range("AB40").Value = -1
$ws.Range("AC40").Value = 0.8999999999999999
$ws.Range("B41").Value = 6853141
$ws.Range("F41").Value = 'Aswan FC'
$ws.Range("G41").Value = 'El Daklyeh'
$ws.Range("H41").Value = 2
$ws.Range("I41").Value = 4
$ws.Range("J41").Value = 'A'
$ws.Range("K41").Value = 2
$ws.Range("L41").Value = 2.75
$ws.Range("M41").Value = 4
$ws.Range("N41").Value = 2.4
$ws.Range("O41").Value = 2.75
$ws.Range("P41").Value = 3
$ws.Range("Q41").Value = -0.25
$ws.Range("R41").Value = 2
$ws.Range("S41").Value = 1.8
$ws.Range("W41").Value = -1
$ws.Range("Y41").Value = 2
$ws.Range("Z41").Value = -1
$ws.Range("AA41").Value = 0.8
$ws.Range("AB41").Value = 0.8999999999999999
$ws.Range("AC41").Value = -1
$ws.Range("B57").Value = 7208756
$ws.Range("F57").Value = 'Al Ittihad Al Sakandary'
$ws.Range("G57").Value = 'Ceramica Cleopatra'
$ws.Range("H57").Value = 1
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 'H'
$ws.Range("K57").Value = 2.8
$ws.Range("L57").Value = 2.9
$ws.Range("M57").Value = 2.55
$ws.Range("N57").Value = 2.8
$ws.Range("O57").Value = 2.9
$ws.Range("P57").Value = 2.55
$ws.Range("Q57").Value = 0
$ws.Range("R57").Value = 1.975
$ws.Range("S57").Value = 1.825
$ws.Range("T57").Value = 2.25
$ws.Range("U57").Value = 2
$ws.Range("V57").Value = 1.8
$ws.Range("W57").Value = 1.8
$ws.Range("X57").Value = -1
$ws.Range("Z57").Value = 0.9750000000000001
$ws.Range("AA57").Value = -1
$ws.Range("AB57").Value = -1
$ws.Range("AC57").Value = 0.8
$ws.Range("B58").Value = 7208367
$ws.Range("F58").Value = 'Pyramids FC'
$ws.Range("G58").Value = 'El Zamalek'
$ws.Range("H58").Value = 2
$ws.Range("I58").Value = 2
$ws.Range("J58").Value = 'D'
$ws.Range("K58").Value = 1.909
$ws.Range("L58").Value = 3.2
$ws.Range("M58").Value = 3.8
$ws.Range("N58").Value = 2.05
$ws.Range("O58").Value = 3.1
$ws.Range("P58").Value = 3.5
$ws.Range("Q58").Value = -0.5
$ws.Range("R58").Value = 2.025
$ws.Range("S58").Value = 1.775
$ws.Range("T58").Value = 2.5
$ws.Range("U58").Value = 1.975
$ws.Range("V58").Value = 1.825
$ws.Range("W58").Value = -1
$ws.Range("X58").Value = 2.1
$ws.Range("Z58").Value = -1
$ws.Range("AA58").Value = 0.7749999999999999
$ws.Range("AB58").Value = 0.9750000000000001
$ws.Range("AC58").Value = -1
$ws.Range("B74").Value = 7217625
$ws.Range("F74").Value = 'Al Ahly Cairo'
$ws.Range("G74").Value = 'Ismaily SC'
$ws.Range("H74").Value = 3
$ws.Range("I74").Value = 1
$ws.Range("K74").Value = 1.25
$ws.Range("L74").Value = 4.75
$ws.Range("M74").Value = 10
$ws.Range("N74").Value = 1.333
$ws.Range("O74").Value = 4.2
$ws.Range("P74").Value = 8
$ws.Range("Q74").Value = -1.25
$ws.Range("R74").Value = 1.775
$ws.Range("S74").Value = 2.025
$ws.Range("U74").Value = 1.875
$ws.Range("V74").Value = 1.925
$ws.Range("W74").Value = 0.333
$ws.Range("Z74").Value = 0.7749999999999999
$ws.Range("AA74").Value = -1
$ws.Range("AB74").Value = 0.875
$ws.Range("AC74").Value = -1
$ws.Range("B75").Value = 7217624
$ws.Range("F75").Value = 'Pyramids FC'
$ws.Range("G75").Value = 'Enppi'
$ws.Range("H75").Value = 1
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 1.444
$ws.Range("L75").Value = 3.75
$ws.Range("M75").Value = 6.5
$ws.Range("N75").Value = 1.5
$ws.Range("O75").Value = 3.75
$ws.Range("P75").Value = 6
$ws.Range("Q75").Value = -1
$ws.Range("R75").Value = 1.85
$ws.Range("S75").Value = 1.95
$ws.Range("U75").Value = 1.975
$ws.Range("V75").Value = 1.825
$ws.Range("W75").Value = 0.5
$ws.Range("Z75").Value = 0
$ws.Range("AA75").Value = -0
$ws.Range("AB75").Value = -1
$ws.Range("AC75").Value = 0.825
$ws.Range("B79").Value = 7217630
$ws.Range("F79").Value = 'El Zamalek'
$ws.Range("G79").Value = 'Smouha'
$ws.Range("H79").Value = 5
$ws.Range("K79").Value = 1.65
$ws.Range("L79").Value = 3.5
$ws.Range("M79").Value = 4.5
$ws.Range("N79").Value = 1.909
$ws.Range("O79").Value = 3.25
$ws.Range("P79").Value = 3.6
$ws.Range("Q79").Value = -0.5
$ws.Range("R79").Value = 1.925
$ws.Range("S79").Value = 1.875
$ws.Range("T79").Value = 2.5
$ws.Range("U79").Value = 1.975
$ws.Range("V79").Value = 1.825
$ws.Range("W79").Value = 0.909
$ws.Range("Z79").Value = 0.925
$ws.Range("AB79").Value = 0.9750000000000001
$ws.Range("B80").Value = 7217629
$ws.Range("F80").Value = 'Enppi'
$ws.Range("G80").Value = 'National Bank'
$ws.Range("H80").Value = 3
$ws.Range("K80").Value = 2.4
$ws.Range("L80").Value = 2.875
$ws.Range("M80").Value = 3
$ws.Range("N80").Value = 2.5
$ws.Range("O80").Value = 2.8
$ws.Range("P80").Value = 2.9
$ws.Range("Q80").Value = 0
$ws.Range("R80").Value = 1.8
$ws.Range("S80").Value = 2
$ws.Range("T80").Value = 2
$ws.Range("U80").Value = 1.85
$ws.Range("V80").Value = 1.95
$ws.Range("W80").Value = 1.5
$ws.Range("Z80").Value = 0.8
$ws.Range("AB80").Value = 0.8500000000000001
$ws.Range("B87").Value = 7217639
$ws.Range("F87").Value = 'ZED FC'
$ws.Range("G87").Value = 'Smouha'
$ws.Range("H87").Value = 1
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 'D'
$ws.Range("K87").Value = 2.1
$ws.Range("L87").Value = 3
$ws.Range("M87").Value = 3.25
$ws.Range("N87").Value = 1.833
$ws.Range("O87").Value = 3.2
$ws.Range("P87").Value = 4
$ws.Range("Q87").Value = -0.5
$ws.Range("R87").Value = 1.9
$ws.Range("S87").Value = 1.9
$ws.Range("T87").Value = 2.25
$ws.Range("U87").Value = 2
$ws.Range("V87").Value = 1.8
$ws.Range("W87").Value = -1
$ws.Range("X87").Value = 2.2
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = 0.8999999999999999
$ws.Range("AB87").Value = -0.5
$ws.Range("AC87").Value = 0.4
$ws.Range("B88").Value = 7217638
$ws.Range("F88").Value = 'Al Ittihad Al Sakandary'
$ws.Range("G88").Value = 'Al Moqawloon Al Arab'
$ws.Range("H88").Value = 3
$ws.Range("I88").Value = 2
$ws.Range("J88").Value = 'H'
$ws.Range("K88").Value = 2.25
$ws.Range("L88").Value = 2.9
$ws.Range("M88").Value = 3.1
$ws.Range("N88").Value = 2.45
$ws.Range("O88").Value = 2.875
$ws.Range("P88").Value = 2.875
$ws.Range("Q88").Value = 0
$ws.Range("R88").Value = 1.725
$ws.Range("S88").Value = 2.075
$ws.Range("T88").Value = 2
$ws.Range("U88").Value = 1.775
$ws.Range("V88").Value = 2.025
$ws.Range("W88").Value = 1.45
$ws.Range("X88").Value = -1
$ws.Range("Z88").Value = 0.7250000000000001
$ws.Range("AA88").Value = -1
$ws.Range("AB88").Value = 0.7749999999999999
$ws.Range("AC88").Value = -1
$ws.Range("B107").Value = 7217659
$ws.Range("F107").Value = 'Pharco FC'
$ws.Range("G107").Value = 'Al Moqawloon Al Arab'
$ws.Range("K107").Value = 2.5
$ws.Range("L107").Value = 2.8
$ws.Range("M107").Value = 2.8
$ws.Range("N107").Value = 2.6
$ws.Range("O107").Value = 2.7
$ws.Range("P107").Value = 2.75
$ws.Range("Q107").Value = 0
$ws.Range("R107").Value = 1.825
$ws.Range("S107").Value = 1.975
$ws.Range("U107").Value = 2
$ws.Range("V107").Value = 1.8
$ws.Range("X107").Value = 1.7
$ws.Range("Z107").Value = 0
$ws.Range("AA107").Value = -0
$ws.Range("B108").Value = 7217658
$ws.Range("F108").Value = 'Talaea El Geish'
$ws.Range("G108").Value = 'El Gounah'
$ws.Range("K108").Value = 2.1
$ws.Range("L108").Value = 2.875
$ws.Range("M108").Value = 3.5
$ws.Range("N108").Value = 2.15
$ws.Range("O108").Value = 2.8
$ws.Range("P108").Value = 3.5
$ws.Range("Q108").Value = -0.25
$ws.Range("R108").Value = 1.85
$ws.Range("S108").Value = 1.95
$ws.Range("U108").Value = 1.825
$ws.Range("V108").Value = 1.975
$ws.Range("X108").Value = 1.8
$ws.Range("Z108").Value = -0.5
$ws.Range("AA108").Value = 0.475
$ws.Range("B172").Value = 7878664
$ws.Range("F172").Value = 'National Bank'
$ws.Range("G172").Value = 'Al Ahly Cairo'
$ws.Range("H172").Value = 4
$ws.Range("I172").Value = 3
$ws.Range("J172").Value = 'H'
$ws.Range("K172").Value = 5
$ws.Range("L172").Value = 4
$ws.Range("M172").Value = 1.5
$ws.Range("N172").Value = 5.5
$ws.Range("O172").Value = 3.4
$ws.Range("P172").Value = 1.55
$ws.Range("Q172").Value = 1
$ws.Range("R172").Value = 1.775
$ws.Range("S172").Value = 2.025
$ws.Range("T172").Value = 2.25
$ws.Range("U172").Value = 1.8
$ws.Range("V172").Value = 2
$ws.Range("W172").Value = 4.5
$ws.Range("X172").Value = -1
$ws.Range("Z172").Value = 0.7749999999999999
$ws.Range("AA172").Value = -1
$ws.Range("AB172").Value = 0.8
$ws.Range("B174").Value = 7946280
$ws.Range("F174").Value = 'El Masry'
$ws.Range("G174").Value = 'ZED FC'
$ws.Range("H174").Value = 2
$ws.Range("I174").Value = 2
$ws.Range("J174").Value = 'D'
$ws.Range("K174").Value = 2.4
$ws.Range("L174").Value = 2.8
$ws.Range("M174").Value = 3.1
$ws.Range("N174").Value = 2.6
$ws.Range("O174").Value = 2.5
$ws.Range("P174").Value = 3.1
$ws.Range("Q174").Value = 0
$ws.Range("R174").Value = 1.725
$ws.Range("S174").Value = 2.075
$ws.Range("T174").Value = 2
$ws.Range("U174").Value = 1.875
$ws.Range("V174").Value = 1.925
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 1.5
$ws.Range("Z174").Value = 0
$ws.Range("AA174").Value = -0
$ws.Range("AB174").Value = 0.875
$ws.Range("E175").Value = 45385.66666666666
$ws.Range("N175").Value = 1.833
$ws.Range("O175").Value = 3.2
$ws.Range("B176").Value = 7217745
$ws.Range("E176").Value = 45385.66666666666
$ws.Range("F176").Value = 'Ismaily SC'
$ws.Range("G176").Value = 'Smouha'
$ws.Range("K176").Value = 2.7
$ws.Range("M176").Value = 2.6
$ws.Range("N176").Value = 2.7
$ws.Range("P176").Value = 2.6
$ws.Range("R176").Value = 2
$ws.Range("S176").Value = 1.8
$ws.Range("U176").Value = 1.825
$ws.Range("V176").Value = 1.975
$ws.Range("B177").Value = 7217744
$ws.Range("E177").Value = 45385.66666666666
$ws.Range("F177").Value = 'El Gounah'
$ws.Range("G177").Value = 'Al Moqawloon Al Arab'
$ws.Range("K177").Value = 2.45
$ws.Range("M177").Value = 2.9
$ws.Range("N177").Value = 2.45
$ws.Range("P177").Value = 2.9
$ws.Range("R177").Value = 1.775
$ws.Range("S177").Value = 2.025
$ws.Range("U177").Value = 1.8
$ws.Range("V177").Value = 2
